$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is bumped by one day
# (45188 -> 45189, i.e. 2023-09-19 -> 2023-09-20) for every data row.
for ($r = 2; $r -le 339; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2() -eq 45188) {
        $cell.Value = 45189
    }
}
